# Updates scraped profit-tracking values (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) for the affected Leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: H9,I9,J9,K9,L9,M9,N9
$ws.Range("H9").Value = 363.33334
$ws.Range("I9").Value = 426
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 426
$ws.Range("L9").Value = 50
$ws.Range("M9").Value = -257
$ws.Range("N9").Value = -388
# Row 28: H28,I28,K28,M28
$ws.Range("H28").Value = 3144.125
$ws.Range("I28").Value = 3021.8572
$ws.Range("K28").Value = 3021.8572
$ws.Range("M28").Value = -2536.8572
# Row 80: H80,I80,K80,M80
$ws.Range("H80").Value = 618.63635
$ws.Range("I80").Value = 525
$ws.Range("K80").Value = 1575
$ws.Range("M80").Value = -577
# Row 83: H83,I83,K83,M83
$ws.Range("H83").Value = 618.63635
$ws.Range("I83").Value = 525
$ws.Range("K83").Value = 4725
$ws.Range("M83").Value = 267
# Row 88: H88,J88,L88,N88
$ws.Range("H88").Value = 9159
$ws.Range("J88").Value = 9568.25
$ws.Range("L88").Value = 9568.25
$ws.Range("N88").Value = -10380.25
# Row 91: H91,J91,L91,N91
$ws.Range("H91").Value = 9159
$ws.Range("J91").Value = 9568.25
$ws.Range("L91").Value = 9568.25
$ws.Range("N91").Value = -12376.25
# Row 132: H132,J132,L132,N132
$ws.Range("H132").Value = 8229
$ws.Range("J132").Value = 9079.799999999999
$ws.Range("L132").Value = 27239.4
$ws.Range("N132").Value = -32299.4
# Row 135: H135,I135,K135,M135
$ws.Range("H135").Value = 1064.7142
$ws.Range("I135").Value = 714.75
$ws.Range("K135").Value = 6432.75
$ws.Range("M135").Value = -3897.75

$ws = $wb.Worksheets.Item("ARM")
# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 3357.111
$ws.Range("I122").Value = 3026.75
$ws.Range("K122").Value = 9080.25
$ws.Range("M122").Value = -6630.25

$ws = $wb.Worksheets.Item("BSM")
# Row 20: H20,I20,K20,M20
$ws.Range("H20").Value = 3336.111
$ws.Range("I20").Value = 3336.111
$ws.Range("K20").Value = 3336.111
$ws.Range("M20").Value = -3089.111
# Row 94: H94,I94,K94,M94
$ws.Range("H94").Value = 1767.3636
$ws.Range("I94").Value = 1737.8889
$ws.Range("K94").Value = 1737.8889
$ws.Range("M94").Value = -1286.8889
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 6815.077
$ws.Range("I134").Value = 3288.7778
$ws.Range("K134").Value = 9866.3334
$ws.Range("M134").Value = -7331.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 16: H16,J16,L16,N16
$ws.Range("H16").Value = 729.5
$ws.Range("J16").Value = 459
$ws.Range("L16").Value = 459
$ws.Range("N16").Value = -1033
# Row 58: H58,J58,L58
$ws.Range("H58").Value = 2443.75
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()  # removed in target row 58
# Row 105: H105,I105,K105,M105
$ws.Range("H105").Value = 2524.4
$ws.Range("I105").Value = 634.8570999999999
$ws.Range("K105").Value = 634.8570999999999
$ws.Range("M105").Value = 1112.1429
# Row 113: H113,J113,L113,N113
$ws.Range("H113").Value = 729.5
$ws.Range("J113").Value = 459
$ws.Range("L113").Value = 459
$ws.Range("N113").Value = -4799
# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 1871.2
$ws.Range("I122").Value = 1880.2222
$ws.Range("K122").Value = 5640.6666
$ws.Range("M122").Value = -3190.6666
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 2277.7778
$ws.Range("I132").Value = 2277.7778
$ws.Range("K132").Value = 6833.3334
$ws.Range("M132").Value = -4303.3334
# Row 136: H136,J136,L136
$ws.Range("H136").Value = 2443.75
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()  # removed in target row 136

$ws = $wb.Worksheets.Item("CUL")
# Row 68: H68,I68,J68,K68,L68,M68,N68
$ws.Range("H68").Value = 874.75
$ws.Range("I68").Value = 899
$ws.Range("J68").Value = 866.6667
$ws.Range("K68").Value = 2697
$ws.Range("L68").Value = 2600.0001
$ws.Range("M68").Value = -1886
$ws.Range("N68").Value = -4222.0001
# Row 71: H71,I71,J71,K71,L71,M71,N71
$ws.Range("H71").Value = 874.75
$ws.Range("I71").Value = 899
$ws.Range("J71").Value = 866.6667
$ws.Range("K71").Value = 8091
$ws.Range("L71").Value = 7800.0003
$ws.Range("M71").Value = -4035
$ws.Range("N71").Value = -15912.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80,J80,L80,N80
$ws.Range("H80").Value = 2241.125
$ws.Range("J80").Value = 1999.5
$ws.Range("L80").Value = 1999.5
$ws.Range("N80").Value = -3995.5
# Row 83: H83,J83,L83,N83
$ws.Range("H83").Value = 2241.125
$ws.Range("J83").Value = 1999.5
$ws.Range("L83").Value = 9997.5
$ws.Range("N83").Value = -19981.5
# Row 97: I97,J97,K97,L97,M97,N97
$ws.Range("I97").Value = 1434.5
$ws.Range("J97").Value = 2199.5
$ws.Range("K97").Value = 1434.5
$ws.Range("L97").Value = 2199.5
$ws.Range("M97").Value = -938.5
$ws.Range("N97").Value = -3191.5
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 3883
$ws.Range("I126").Value = 4276.2
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 12828.6
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -10358.6
$ws.Range("N126").Value = -13640

$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7,I7,K7,M7
$ws.Range("H7").Value = 4575.375
$ws.Range("I7").Value = 4086.1428
$ws.Range("K7").Value = 4086.1428
$ws.Range("M7").Value = -3974.1428
# Row 55: H55,I55,J55,K55,L55,M55,N55
$ws.Range("H55").Value = 2654.2222
$ws.Range("I55").Value = 3347.25
$ws.Range("J55").Value = 2099.8
$ws.Range("K55").Value = 3347.25
$ws.Range("L55").Value = 2099.8
$ws.Range("M55").Value = -3174.25
$ws.Range("N55").Value = -2445.8
# Row 124: H124,J124,L124,N124
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 4575.375
$ws.Range("I126").Value = 4086.1428
$ws.Range("K126").Value = 12258.4284
$ws.Range("M126").Value = -9788.428400000001
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 7378.8667
$ws.Range("I136").Value = 3744.818
$ws.Range("K136").Value = 11234.454
$ws.Range("M136").Value = -8684.454000000002

$ws = $wb.Worksheets.Item("WVR")
# Row 81: H81,I81,J81,K81,L81,M81,N81
$ws.Range("H81").Value = 7741.4165
$ws.Range("I81").Value = 8649.700000000001
$ws.Range("J81").Value = 3200
$ws.Range("K81").Value = 17299.4
$ws.Range("L81").Value = 6400
$ws.Range("M81").Value = -16238.4
$ws.Range("N81").Value = -8522
# Row 84: H84,I84,J84,K84,L84,M84,N84
$ws.Range("H84").Value = 7741.4165
$ws.Range("I84").Value = 8649.700000000001
$ws.Range("J84").Value = 3200
$ws.Range("K84").Value = 86497
$ws.Range("L84").Value = 32000
$ws.Range("M84").Value = -81193
$ws.Range("N84").Value = -42608
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 578.82355
$ws.Range("I107").Value = 445.16666
$ws.Range("J107").Value = 899.6
$ws.Range("K107").Value = 1335.49998
$ws.Range("L107").Value = 2698.8
$ws.Range("M107").Value = 584.5000199999999
$ws.Range("N107").Value = -6538.8
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 7226.65
$ws.Range("I132").Value = 3666.7273
$ws.Range("K132").Value = 11000.1819
$ws.Range("M132").Value = -8470.1819
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 2047.1666
$ws.Range("I136").Value = 2047.1666
$ws.Range("K136").Value = 6141.4998
$ws.Range("M136").Value = -3591.4998
